$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 640.1539
$ws.Range("I11").Value = 640.1539
$ws.Range("K11").Value = 640.1539
$ws.Range("M11").Value = -500.1539
$ws.Range("H51").Value = 4707
$ws.Range("J51").Value = 5060.5
$ws.Range("L51").Value = 5060.5
$ws.Range("N51").Value = -6028.5
$ws.Range("H116").Value = 11573.85
$ws.Range("J116").Value = 11306.846
$ws.Range("L116").Value = 11306.846
$ws.Range("N116").Value = -18190.846
$ws.Range("H137").Value = 10827.353
$ws.Range("I137").Value = 14577.454
$ws.Range("J137").Value = 3952.1667
$ws.Range("K137").Value = 43732.362
$ws.Range("L137").Value = 11856.5001
$ws.Range("M137").Value = -41182.362
$ws.Range("N137").Value = -16956.5001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4304.3657
$ws.Range("I32").Value = 4386.6455
$ws.Range("J32").Value = 2137.6667
$ws.Range("K32").Value = 4386.6455
$ws.Range("L32").Value = 2137.6667
$ws.Range("M32").Value = -4099.6455
$ws.Range("N32").Value = -2711.6667

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 15874269
$ws.Range("J64").Value = 775
$ws.Range("L64").Value = 775
$ws.Range("N64").Value = -1225
$ws.Range("H67").Value = 15874269
$ws.Range("J67").Value = 775
$ws.Range("L67").Value = 775
$ws.Range("N67").Value = -2335
$ws.Range("H94").Value = 407.06897
$ws.Range("I94").Value = 324.54544
$ws.Range("J94").Value = 666.4286
$ws.Range("K94").Value = 324.54544
$ws.Range("L94").Value = 666.4286
$ws.Range("M94").Value = 126.45456
$ws.Range("N94").Value = -1568.4286

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2414.3547
$ws.Range("I31").Value = 1829.5
$ws.Range("J31").Value = 3477.7273
$ws.Range("K31").Value = 1829.5
$ws.Range("L31").Value = 3477.7273
$ws.Range("M31").Value = -1534.5
$ws.Range("N31").Value = -4067.7273
$ws.Range("H34").Value = 2414.3547
$ws.Range("I34").Value = 1829.5
$ws.Range("J34").Value = 3477.7273
$ws.Range("K34").Value = 1829.5
$ws.Range("L34").Value = 3477.7273
$ws.Range("M34").Value = -1627.5
$ws.Range("N34").Value = -3881.7273
$ws.Range("H58").Value = 1649.9111
$ws.Range("I58").Value = 1528.6666
$ws.Range("J58").Value = 2438
$ws.Range("K58").Value = 1528.6666
$ws.Range("L58").Value = 2438
$ws.Range("M58").Value = -1325.6666
$ws.Range("N58").Value = -2844
$ws.Range("H99").Value = 12128.04
$ws.Range("I99").Value = 7442.643
$ws.Range("K99").Value = 7442.643
$ws.Range("M99").Value = -5944.643
$ws.Range("H126").Value = 12128.04
$ws.Range("I126").Value = 7442.643
$ws.Range("K126").Value = 22327.929
$ws.Range("M126").Value = -19857.929
$ws.Range("H132").Value = 6854.246
$ws.Range("I132").Value = 3555.7659
$ws.Range("K132").Value = 10667.2977
$ws.Range("M132").Value = -8137.297699999999
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H135").Value = 144997.33
$ws.Range("J135").Value = 144997.33
$ws.Range("L135").Value = 144997.33
$ws.Range("N135").Value = -155137.33
$ws.Range("H136").Value = 1649.9111
$ws.Range("I136").Value = 1528.6666
$ws.Range("J136").Value = 2438
$ws.Range("K136").Value = 4585.9998
$ws.Range("L136").Value = 7314
$ws.Range("M136").Value = -2035.9998
$ws.Range("N136").Value = -12414
$ws.Range("H141").Value = 146480.75
$ws.Range("J141").Value = 159232.5
$ws.Range("L141").Value = 159232.5
$ws.Range("N141").Value = -169592.5
$ws.Range("N133").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 384.7143
$ws.Range("I14").Value = 384.7143
$ws.Range("K14").Value = 1154.1429
$ws.Range("M14").Value = -981.1428999999998
$ws.Range("H132").Value = 2622.9375
$ws.Range("I132").Value = 957.8
$ws.Range("J132").Value = 3379.818
$ws.Range("K132").Value = 8620.199999999999
$ws.Range("L132").Value = 30418.362
$ws.Range("M132").Value = -6090.199999999999
$ws.Range("N132").Value = -35478.362
$ws.Range("H138").Value = 3584
$ws.Range("I138").Value = 3584
$ws.Range("K138").Value = 10752
$ws.Range("M138").Value = -5612
$ws.Range("H139").Value = 2956.2
$ws.Range("I139").Value = 1487.9231
$ws.Range("J139").Value = 12500
$ws.Range("K139").Value = 4463.7693
$ws.Range("L139").Value = 37500
$ws.Range("M139").Value = 676.2307000000001
$ws.Range("N139").Value = -47780

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2265.2424
$ws.Range("I113").Value = 1905.5238
$ws.Range("J113").Value = 2894.75
$ws.Range("K113").Value = 1905.5238
$ws.Range("L113").Value = 2894.75
$ws.Range("M113").Value = 264.4762000000001
$ws.Range("N113").Value = -7234.75
$ws.Range("H132").Value = 2449.5
$ws.Range("I132").Value = 1890.2041
$ws.Range("J132").Value = 3891.8948
$ws.Range("K132").Value = 5670.6123
$ws.Range("L132").Value = 11675.6844
$ws.Range("M132").Value = -3140.6123
$ws.Range("N132").Value = -16735.6844

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 995
$ws.Range("I22").Value = 995
$ws.Range("J22").Value = 995
$ws.Range("K22").Value = 995
$ws.Range("L22").Value = 995
$ws.Range("M22").Value = -700
$ws.Range("N22").Value = -1585
$ws.Range("H27").Value = 995
$ws.Range("I27").Value = 995
$ws.Range("J27").Value = 995
$ws.Range("K27").Value = 995
$ws.Range("L27").Value = 995
$ws.Range("M27").Value = -888
$ws.Range("N27").Value = -1209
$ws.Range("H132").Value = 45676.945
$ws.Range("I132").Value = 50804.152
$ws.Range("J132").Value = 3377.5
$ws.Range("K132").Value = 152412.456
$ws.Range("L132").Value = 10132.5
$ws.Range("M132").Value = -149882.456
$ws.Range("N132").Value = -15192.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 73250
$ws.Range("I2").Value = 100000
$ws.Range("J2").Value = 46500
$ws.Range("K2").Value = 100000
$ws.Range("L2").Value = 46500
$ws.Range("M2").Value = -99888
$ws.Range("N2").Value = -46724
$ws.Range("H132").Value = 3091.25
$ws.Range("I132").Value = 2833.6924
$ws.Range("J132").Value = 3760.9
$ws.Range("K132").Value = 8501.0772
$ws.Range("L132").Value = 11282.7
$ws.Range("M132").Value = -5971.0772
$ws.Range("N132").Value = -16342.7
$ws.Range("H136").Value = 2310.394
$ws.Range("I136").Value = 2310.394
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6931.181999999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4381.181999999999
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("N137").ClearContents()
